$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handoff
#   - "29bc900b-..." file is now in translation (handoff just regenerated)
#   - "bec7f7e5-..." file is ready for handoff, but its handback is stale
#     (a newer source revision exists), so an error detail is attached.
# ---------------------------------------------------------------------------

# --- Overview sheet ---------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2016-10-19 11:47:16"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-19 11:47:16"

# Narrow the zh-cn / de-de status columns (they no longer need to be as wide).
$wsOverview.Range("E1").ColumnWidth = 16.333333333333336
$wsOverview.Range("F1").ColumnWidth = 16.333333333333336

# --- zh-cn sheet -------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("H2").Value = "2016-10-19 11:47:06"
$wsZhCn.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11772e054ad11acbf9733480d64c935eba11add6/e2e/29bc900b-f4ef-4d73-b187-8d5b003fec25.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe2efc69e266442dbf1d5ff54e7b04d1f4e6f2e3/e2e/29bc900b-f4ef-4d73-b187-8d5b003fec25.md."

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-10-19 11:47:06"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11772e054ad11acbf9733480d64c935eba11add6/e2e/bec7f7e5-7c4d-4976-b03e-878d90c2a8d6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe2efc69e266442dbf1d5ff54e7b04d1f4e6f2e3/e2e/bec7f7e5-7c4d-4976-b03e-878d90c2a8d6.md."

$wsZhCn.Range("C1").ColumnWidth = 16.333333333333336
$wsZhCn.Range("P1").ColumnWidth = 39.16666666666667

# --- de-de sheet ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("H2").Value = "2016-10-19 11:47:16"
$wsDeDe.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11772e054ad11acbf9733480d64c935eba11add6/e2e/29bc900b-f4ef-4d73-b187-8d5b003fec25.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe2efc69e266442dbf1d5ff54e7b04d1f4e6f2e3/e2e/29bc900b-f4ef-4d73-b187-8d5b003fec25.md."

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-10-19 11:47:16"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11772e054ad11acbf9733480d64c935eba11add6/e2e/bec7f7e5-7c4d-4976-b03e-878d90c2a8d6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe2efc69e266442dbf1d5ff54e7b04d1f4e6f2e3/e2e/bec7f7e5-7c4d-4976-b03e-878d90c2a8d6.md."

$wsDeDe.Range("C1").ColumnWidth = 16.333333333333336
$wsDeDe.Range("P1").ColumnWidth = 39.16666666666667
